$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45202 -> 45203) for every data row (rows 2 through 503).
$lastRow = $ws.Range("A1").End(4).Row
if ($lastRow -lt 503) { $lastRow = 503 }

$rng = $ws.Range("C2:C$lastRow")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2() -eq 45202) {
        $cell.Value = 45203
    }
}
